# Add new rows to Feuil1 describing method_type factor aesthetics
# (label/colour/order used for plotting figure 3 as a stacked barplot
# of method types as a proportion of all predicted articles per oro type)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "method_type"
$ws.Range("B12").Value = "Empirical"
$ws.Range("C12").Value = "Empirical"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "#1b9e77"

$ws.Range("A13").Value = "method_type"
$ws.Range("B13").Value = "Social_primary"
$ws.Range("C13").Value = "Social primary"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "#d95f02"

$ws.Range("A14").Value = "method_type"
$ws.Range("B14").Value = "Mathematical_predictionsimulation"
$ws.Range("C14").Value = "Prediction/simulation"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "#7570b3"

$ws.Range("A15").Value = "method_type"
$ws.Range("B15").Value = "Other"
$ws.Range("C15").Value = "Other"
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = "darkgrey"

# Update the active selection to mirror the authored workbook state
$ws.Range("E16").Select()
